# Generate Report for Handoff
# Update the "bec9db9b-4814-4cc0-a64f-6ea3ce7bf266.md" row (row 3) across the
# Overview / zh-cn / de-de sheets to reflect that the file is now
# "Ready for handoff" with a fresh handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-28-11 08:28:45"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-11 08:28:42"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-11 08:28:45"
